$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "not considered"
$ws.Range("E2").Value = "not considered"
$ws.Range("F2").Value = "not considered"
$ws.Range("G2").Value = "important"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "not considered"
$ws.Range("J2").Value = "considered"
